$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded; it belongs chronologically
# right before the existing row 158, so insert a fresh row there and push
# every following row (old 158-260) down by one (new 159-261).
$ws.Rows("158:158").Insert()

# Populate the newly inserted row 158 with the new data point.
$ws.Range("A158").Value = 10
$ws.Range("B158").Value = "Vega Modelo de Temuco"
$ws.Range("C158").Value = "La Araucanía"
$ws.Range("D158").Value = 44438
$ws.Range("E158").Value = 9
$ws.Range("F158").Value = "Fruta"
$ws.Range("G158").Value = 100108
$ws.Range("H158").Value = "Tropicales y subtropicales"
$ws.Range("I158").Value = 100108005
$ws.Range("J158").Value = "Piña"
$ws.Range("K158").Value = "Caramelo"
$ws.Range("L158").Value = "Segunda"
$ws.Range("M158").Value = 120
$ws.Range("N158").Value = 20000
$ws.Range("O158").Value = 21000
$ws.Range("P158").Value = 20458
$ws.Range("Q158").Value = "$/caja 14 unidades"
$ws.Range("R158").Value = "Ecuador"
$ws.Range("S158").Value = 1461
$ws.Range("T158").Value = 14
